$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "19 loka"
$ws.Range("B10").Value = "18.30-21.00"
$ws.Range("C10").Value = "Kameraluokka, liikkuminen scenessä ja ruudukko"
$ws.Range("D10").Value = "Alussa vaikutti olevan hyvä meininki mutta jumiin jäätiin. Näyttää olevan kantava teema."
$ws.Range("E10").Value = "Kyllähän tämä tästä, hieman siistimmin taas muutama asia opittu kuten compile time constant eli constexpr jota käytettiin kameran alustuksessa."
$ws.Range("G10").Value = 2.5

$ws.Range("A11").Value = "20 loka"
$ws.Range("B11").Value = "20.30-21.30"

$ws.Range("D10").WrapText = $true
$ws.Range("E10").WrapText = $true

$ws.Rows.Item(10).RowHeight = 58

$ws.Range("B11").NumberFormat = $ws.Range("B9").NumberFormat
